$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.592.44"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "'3.513.33"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'586.53"
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("D6").Value = "'132.90"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").Value = "'3.514.93"
$ws.Range("E7").Value = "  -1.09%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("D11").Value = "'7.16"
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").Value = "'0.389"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "'4.105.94"
$ws.Range("E13").Value = "  -1.25%  "
$ws.Range("D14").Value = "'27.88"
$ws.Range("E14").Value = "  +3.66%  "
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("D17").Value = "'3.508.02"
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("D18").Value = "'64.617.97"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("D20").Value = "'14.29"
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("E21").Value = "  -2.22%  "
$ws.Range("D22").Value = "'393.47"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("D23").Value = "'0.578"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").Value = "'3.653.48"
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").Value = "'74.30"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").Value = "'0.0000111"
$ws.Range("E27").Value = "  -2.71%  "
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("D29").Value = "'7.49"
$ws.Range("E29").Value = "  -4.08%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "'2.27"
$ws.Range("E31").Value = "  -0.87%  "
$ws.Range("E32").Value = "  -4.00%  "
$ws.Range("D33").Value = "'3.518.19"
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").Value = "'24.03"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("E37").Value = "  +4.50%  "
$ws.Range("D38").Value = "'5.27"
$ws.Range("E38").Value = "  +4.45%  "
$ws.Range("D39").Value = "'171.27"
$ws.Range("E40").Value = "  +0.89%  "
$ws.Range("D41").Value = "'0.0815"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D42").Value = "'0.814"
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("D43").Value = "'26.49"
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "'42.28"
$ws.Range("E45").Value = "  -1.61%  "
$ws.Range("E46").Value = "  -3.33%  "
$ws.Range("D47").Value = "'4.42"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("D49").Value = "'2.465.85"
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("D50").Value = "'6.89"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").Value = "'0.911"
